$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily log entry for 2026/01/27 is inserted right before the old
# row 707 (2026/12/29 火 13 9). Everything from the old row 707 onward
# shifts down by one row; their values are untouched by the insert.
$ws.Rows(707).Insert()

# Column A holds dates as plain text (e.g. "2026/12/29"), never real
# Excel date serials. Force text formatting before typing the value so
# Excel doesn't auto-convert "2026/01/27" into a date number, then clear
# the formatting again so the new cell ends up styleless like its
# neighbours.
$ws.Range("A707").NumberFormat = "@"
$ws.Range("A707").Value = "2026/01/27"
$ws.Range("B707").Value = "火"
$ws.Range("C707").Value = 15
$ws.Range("D707").Value = 61
$ws.Range("A707:D707").ClearFormats()
